$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 70: "901. Online Stock Span" (Stacks, Medium) ---
# Fill the row in the same order the original author appears to have used:
# Question, Difficulty, Pattern, Link, then Notes (this keeps the new
# shared-string insertion order lined up with the source workbook).

$ws.Range("A70").Value = "901. Online Stock Span"

$ws.Range("B70").Value = "Medium"
$ws.Range("B70").Interior.Color = 49407

$ws.Range("C70").Value = "Stacks"

$link = "https://leetcode.com/problems/online-stock-span/solutions/168311/c-java-python-o-1/?envType=study-plan-v2&envId=leetcode-75 "
$ws.Range("E70").Value = $link
$ws.Hyperlinks.Add($ws.Range("E70"), $link)
$ws.Range("E70").Style = "Hyperlink"

$ws.Range("D70").Value = "Push every pair of <price, result> to a stack. Pop the lower price from the stack and accumulate the count. We can keep popping from the stack when we call next() because we store the result at the time of pushing each element. Thus we maintain the stack invariant and we effectively collapse the spans into a single value."

# --- Update the saved selection to reflect where the author left off ---
$ws.Range("D76").Select() | Out-Null
